$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''23.471.27'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -0.89%  '

$ws.Range("D3").Value = '''1.639.61'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -0.67%  '

$ws.Range("D4").Value = '''1.000'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.13%  '

$ws.Range("D5").Value = '''0.9994'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.16%  '

$ws.Range("D6").Value = '''304.19'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.82%  '

$ws.Range("D7").Value = '''0.3791'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +0.39%  '

$ws.Range("D8").Value = '''51.65'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -2.24%  '

$ws.Range("D9").Value = '''0.3630'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -1.15%  '

$ws.Range("D10").Value = '''0.08198'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +0.43%  '

$ws.Range("D11").Value = '''1.235'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -2.79%  '

$ws.Range("D12").Value = '''0.9981'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -0.38%  '

$ws.Range("D13").Value = '''22.55'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -2.39%  '

$ws.Range("D14").Value = '''6.467'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -3.71%  '

$ws.Range("E15").Value = '  -0.10%  '

$ws.Range("D16").Value = '''0.00001242'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -2.46%  '

$ws.Range("D17").Value = '''1.634.83'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -1.12%  '

$ws.Range("D18").Value = '''95.44'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +0.20%  '

$ws.Range("D19").Value = '''0.06932'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -0.09%  '

$ws.Range("D20").Value = '''6.584'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.46%  '

$ws.Range("D21").Value = '''17.53'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -4.78%  '

$ws.Range("D22").Value = '''0.9995'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.05%  '

$ws.Range("D23").Value = '''12.57'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -3.29%  '

$ws.Range("D24").Value = '''23.480.21'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -0.87%  '

$ws.Range("D25").Value = '''2.505'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +3.26%  '

$ws.Range("D26").Value = '''3.061'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -5.74%  '

$ws.Range("D27").Value = '''21.15'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -1.23%  '

$ws.Range("D28").Value = '''151.42'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -0.37%  '

$ws.Range("E29").Value = '  -0.63%  '

$ws.Range("D30").Value = '''133.59'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -2.54%  '

$ws.Range("D31").Value = '''1.817.99'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -0.98%  '

$ws.Range("D32").Value = '''2.182'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -5.63%  '

$ws.Range("D33").Value = '''6.647'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -6.56%  '

$ws.Range("D34").Value = '''1.064'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +8.88%  '

$ws.Range("D35").Value = '''11.41'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +3.43%  '

$ws.Range("D36").Value = '''0.02765'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -4.40%  '

$ws.Range("E37").Value = '  -3.47%  '

$ws.Range("D38").Value = '''0.08775'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -1.30%  '

$ws.Range("D39").Value = '''0.07129'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -3.04%  '

$ws.Range("D40").Value = '''6.017'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -5.64%  '

$ws.Range("D41").Value = '''0.7067'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -2.24%  '

$ws.Range("D42").Value = '''1.339'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -3.61%  '

$ws.Range("B43").Value = 'EnergySwap'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D43").Value = '''15.82'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -4.19%  '

$ws.Range("B44").Value = 'Aptos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D44").Value = '''12.17'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -4.65%  '

$ws.Range("D45").Value = '''0.6557'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -1.62%  '

$ws.Range("D46").Value = '''0.9991'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +0.10%  '

$ws.Range("D47").Value = '''2.286'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -4.58%  '

$ws.Range("D48").Value = '''3.967'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -1.56%  '

$ws.Range("D49").Value = '''0.07985'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -1.16%  '

$ws.Range("D50").Value = '''127.67'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -1.07%  '

$ws.Range("D51").Value = '''1.196'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -2.76%  '
